$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 80553
$ws.Range("E2").Value = 1413
$ws.Range("F2").Value = 1413
$ws.Range("G2").Value = 509
$ws.Range("H2").Value = 114
$ws.Range("I2").Value = 567
$ws.Range("J2").Value = -454
$ws.Range("K2").Value = 125970
$ws.Range("L2").Value = 78175
$ws.Range("M2").Value = 47795
$ws.Range("N2").Value = 43741
$ws.Range("O2").Value = 4054
$ws.Range("P2").Value = 8153
$ws.Range("Q2").Value = 3625
$ws.Range("R2").Value = -2850
$ws.Range("S2").Value = -4517
$ws.Range("T2").Value = 4657
$ws.Range("U2").Value = -1032
$ws.Range("V2").Value = 49611
$ws.Range("W2").Value = 1.75
$ws.Range("X2").Value = 0.14
$ws.Range("Y2").Value = 1.35
$ws.Range("Z2").Value = 0.09
$ws.Range("AA2").Value = 163.56
$ws.Range("AB2").Value = 423.44
$ws.Range("AC2").Value = 364
$ws.Range("AD2").Value = 32.45
$ws.Range("AE2").Value = 26824
$ws.Range("AF2").Value = 0.44
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.27
$ws.Range("AI2").Value = 43.21
$ws.Range("AJ2").Value = 161939571

# Row 3
$ws.Range("D3").Value = 80370
$ws.Range("E3").Value = 3370
$ws.Range("F3").Value = 3370
$ws.Range("G3").Value = 2284
$ws.Range("H3").Value = 1804
$ws.Range("I3").Value = 1882
$ws.Range("J3").Value = -77
$ws.Range("K3").Value = 138526
$ws.Range("L3").Value = 89215
$ws.Range("M3").Value = 49312
$ws.Range("N3").Value = 46334
$ws.Range("O3").Value = 2978
$ws.Range("P3").Value = 8153
$ws.Range("Q3").Value = 12093
$ws.Range("R3").Value = -9526
$ws.Range("S3").Value = -1540
$ws.Range("T3").Value = 7592
$ws.Range("U3").Value = 4501
$ws.Range("V3").Value = 50411
$ws.Range("W3").Value = 4.19
$ws.Range("X3").Value = 2.25
$ws.Range("Y3").Value = 4.18
$ws.Range("Z3").Value = 1.36
$ws.Range("AA3").Value = 180.92
$ws.Range("AB3").Value = 455.82
$ws.Range("AC3").Value = 1154
$ws.Range("AD3").Value = 23.57
$ws.Range("AE3").Value = 28415
$ws.Range("AF3").Value = 0.96
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 0.55
$ws.Range("AI3").Value = 13.03
$ws.Range("AJ3").Value = 161939571

# Row 4
$ws.Range("D4").Value = 92588
$ws.Range("E4").Value = 7792
$ws.Range("F4").Value = 7792
$ws.Range("G4").Value = 8864
$ws.Range("H4").Value = 7709
$ws.Range("I4").Value = 7553
$ws.Range("J4").Value = 157
$ws.Range("K4").Value = 138179
$ws.Range("L4").Value = 83496
$ws.Range("M4").Value = 54683
$ws.Range("N4").Value = 52849
$ws.Range("O4").Value = 1834
$ws.Range("P4").Value = 8297
$ws.Range("Q4").Value = 11510
$ws.Range("R4").Value = -4601
$ws.Range("S4").Value = -3182
$ws.Range("T4").Value = 6148
$ws.Range("U4").Value = 5362
$ws.Range("V4").Value = 48296
$ws.Range("W4").Value = 8.42
$ws.Range("X4").Value = 8.33
$ws.Range("Y4").Value = 15.23
$ws.Range("Z4").Value = 5.57
$ws.Range("AA4").Value = 152.69
$ws.Range("AB4").Value = 532.12
$ws.Range("AC4").Value = 4564
$ws.Range("AD4").Value = 5.41
$ws.Range("AE4").Value = 32179
$ws.Range("AF4").Value = 0.77
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 1.42
$ws.Range("AI4").Value = 7.62
$ws.Range("AJ4").Value = 164809359

# Row 5
$ws.Range("D5").Value = 93418
$ws.Range("E5").Value = 7564
$ws.Range("F5").Value = 7564
$ws.Range("G5").Value = 10659
$ws.Range("H5").Value = 8345
$ws.Range("I5").Value = 8247
$ws.Range("J5").Value = 98
$ws.Range("K5").Value = 136495
$ws.Range("L5").Value = 74620
$ws.Range("M5").Value = 61875
$ws.Range("N5").Value = 60151
$ws.Range("O5").Value = 1723
$ws.Range("P5").Value = 8297
$ws.Range("Q5").Value = 9162
$ws.Range("R5").Value = -4032
$ws.Range("S5").Value = -6856
$ws.Range("T5").Value = 4924
$ws.Range("U5").Value = 4239
$ws.Range("V5").Value = 44174
$ws.Range("W5").Value = 8.1
$ws.Range("X5").Value = 8.93
$ws.Range("Y5").Value = 14.6
$ws.Range("Z5").Value = 6.08
$ws.Range("AA5").Value = 120.6
$ws.Range("AB5").Value = 626
$ws.Range("AC5").Value = 4970
$ws.Range("AD5").Value = 6.36
$ws.Range("AE5").Value = 36625
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 1.11
$ws.Range("AI5").Value = 6.98
$ws.Range("AJ5").Value = 164809359

# Row 6
$ws.Range("D6").Value = 90460
$ws.Range("E6").Value = 3543
$ws.Range("F6").Value = 3543
$ws.Range("G6").Value = 2883
$ws.Range("H6").Value = 1604
$ws.Range("I6").Value = 1867
$ws.Range("K6").Value = 152315
$ws.Range("L6").Value = 90044
$ws.Range("M6").Value = 62271
$ws.Range("N6").Value = 60900
$ws.Range("P6").Value = 8212
$ws.Range("Q6").Value = 7267
$ws.Range("R6").Value = -4346
$ws.Range("S6").Value = -696
$ws.Range("T6").Value = 7201
$ws.Range("U6").Value = 66
$ws.Range("V6").Value = 58885
$ws.Range("W6").Value = 3.92
$ws.Range("X6").Value = 1.77
$ws.Range("Y6").Value = 3.08
$ws.Range("Z6").Value = 1.11
$ws.Range("AA6").Value = 144.6
$ws.Range("AB6").Value = 651.9299999999999
$ws.Range("AC6").Value = 1133
$ws.Range("AD6").Value = 17.83
$ws.Range("AE6").Value = 37081
$ws.Range("AF6").Value = 0.54
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 0.99
$ws.Range("AI6").Value = 17.63
$ws.Range("AJ6").Value = 163110394

# Row 7
$ws.Range("D7").Value = 95583
$ws.Range("E7").Value = 4340
$ws.Range("G7").Value = 3554
$ws.Range("H7").Value = 2664
$ws.Range("I7").Value = 2685
$ws.Range("K7").Value = 158389
$ws.Range("L7").Value = 94504
$ws.Range("M7").Value = 63885
$ws.Range("N7").Value = 62615
$ws.Range("P7").Value = 8188
$ws.Range("Q7").Value = 7469
$ws.Range("R7").Value = -8706
$ws.Range("S7").Value = -255
$ws.Range("T7").Value = 9471
$ws.Range("U7").Value = -2667
$ws.Range("W7").Value = 4.54
$ws.Range("X7").Value = 2.79
$ws.Range("Y7").Value = 4.35
$ws.Range("Z7").Value = 1.71
$ws.Range("AA7").Value = 147.93
$ws.Range("AC7").Value = 1647
$ws.Range("AD7").Value = 10.23
$ws.Range("AE7").Value = 38508
$ws.Range("AF7").Value = 0.44
$ws.Range("AG7").Value = 246
$ws.Range("AH7").Value = 1.46
$ws.Range("AI7").Value = 14.82

# Row 8
$ws.Range("D8").Value = 102414
$ws.Range("E8").Value = 5254
$ws.Range("G8").Value = 5420
$ws.Range("H8").Value = 4064
$ws.Range("I8").Value = 4028
$ws.Range("K8").Value = 164479
$ws.Range("L8").Value = 97092
$ws.Range("M8").Value = 67387
$ws.Range("N8").Value = 66008
$ws.Range("P8").Value = 8188
$ws.Range("Q8").Value = 7980
$ws.Range("R8").Value = -6828
$ws.Range("S8").Value = -399
$ws.Range("T8").Value = 7017
$ws.Range("U8").Value = 625
$ws.Range("W8").Value = 5.13
$ws.Range("X8").Value = 3.97
$ws.Range("Y8").Value = 6.26
$ws.Range("Z8").Value = 2.52
$ws.Range("AA8").Value = 144.08
$ws.Range("AC8").Value = 2477
$ws.Range("AD8").Value = 6.8
$ws.Range("AE8").Value = 40595
$ws.Range("AF8").Value = 0.42
$ws.Range("AG8").Value = 279
$ws.Range("AH8").Value = 1.65
$ws.Range("AI8").Value = 11.17

# Row 9
$ws.Range("D9").Value = 104513
$ws.Range("E9").Value = 5807
$ws.Range("G9").Value = 6079
$ws.Range("H9").Value = 4545
$ws.Range("I9").Value = 4616
$ws.Range("K9").Value = 169710
$ws.Range("L9").Value = 98401
$ws.Range("M9").Value = 71310
$ws.Range("N9").Value = 69842
$ws.Range("P9").Value = 8188
$ws.Range("Q9").Value = 8903
$ws.Range("R9").Value = -6351
$ws.Range("S9").Value = -949
$ws.Range("T9").Value = 6813
$ws.Range("U9").Value = 1876
$ws.Range("W9").Value = 5.56
$ws.Range("X9").Value = 4.35
$ws.Range("Y9").Value = 6.8
$ws.Range("Z9").Value = 2.72
$ws.Range("AA9").Value = 137.99
$ws.Range("AC9").Value = 2839
$ws.Range("AD9").Value = 5.94
$ws.Range("AE9").Value = 42952
$ws.Range("AF9").Value = 0.39
$ws.Range("AG9").Value = 293
$ws.Range("AH9").Value = 1.74
$ws.Range("AI9").Value = 10.24
